# Update the cryptos price list: refresh Price/Volume(1h) figures for most rows,
# and swap the Fetch.AI/Bittensor and Aave/USDe row pairs (rows 31-32 and 40-41).
# Note: several "Price" values are single decimal numbers (e.g. "1.40", "0.379")
# that Excel's COM layer would otherwise auto-convert to numeric types when
# assigned via .Value. Prefixing with a literal apostrophe forces Excel to
# keep them as text, matching the original inlineStr/text cell type, while
# two-dot "thousands.hundreds.cents" style values (e.g. "75.411.11") are
# already unambiguous text and don't need the prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '75.411.11'
$ws.Range("E2").Value = '  +1.42%  '
$ws.Range("D3").Value = '2.827.65'
$ws.Range("E3").Value = '  +6.20%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''189.84'
$ws.Range("E5").Value = '  +1.68%  '
$ws.Range("D6").Value = '''595.65'
$ws.Range("E6").Value = '  +1.94%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''0.547'
$ws.Range("E8").Value = '  +3.12%  '
$ws.Range("E9").Value = '  -1.37%  '
$ws.Range("D10").Value = '2.825.32'
$ws.Range("E10").Value = '  +6.12%  '
$ws.Range("D11").Value = '''0.379'
$ws.Range("E11").Value = '  +7.03%  '
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("D13").Value = '''4.90'
$ws.Range("E13").Value = '  +4.63%  '
$ws.Range("D14").Value = '3.344.68'
$ws.Range("E14").Value = '  +6.28%  '
$ws.Range("D15").Value = '75.337.55'
$ws.Range("E15").Value = '  +1.31%  '
$ws.Range("D16").Value = '''0.0000188'
$ws.Range("E16").Value = '  +1.22%  '
$ws.Range("D17").Value = '''27.01'
$ws.Range("E17").Value = '  +2.90%  '
$ws.Range("D18").Value = '2.815.96'
$ws.Range("E18").Value = '  +6.12%  '
$ws.Range("D19").Value = '''8.89'
$ws.Range("E19").Value = '  -4.17%  '
$ws.Range("D20").Value = '''12.34'
$ws.Range("E20").Value = '  +4.01%  '
$ws.Range("D21").Value = '''378.61'
$ws.Range("E21").Value = '  +3.78%  '
$ws.Range("D22").Value = '''2.30'
$ws.Range("E22").Value = '  +1.54%  '
$ws.Range("D23").Value = '''4.11'
$ws.Range("E23").Value = '  +1.53%  '
$ws.Range("E24").Value = '  -0.54%  '
$ws.Range("D25").Value = '''71.02'
$ws.Range("E25").Value = '  +1.81%  '
$ws.Range("D26").Value = '2.961.72'
$ws.Range("E26").Value = '  +7.41%  '
$ws.Range("D27").Value = '''4.18'
$ws.Range("E27").Value = '  +2.20%  '
$ws.Range("D28").Value = '''9.77'
$ws.Range("E28").Value = '  +5.64%  '
$ws.Range("D29").Value = '''0.0000104'
$ws.Range("E29").Value = '  +11.93%  '
$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").Value = '''1.40'
$ws.Range("E31").Value = '  +1.40%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").Value = '''516.10'
$ws.Range("E32").Value = '  -0.57%  '
$ws.Range("D33").Value = '''7.69'
$ws.Range("E33").Value = '  +0.94%  '
$ws.Range("D34").Value = '''1.81'
$ws.Range("E34").Value = '  +4.17%  '
$ws.Range("D35").Value = '''0.999'
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").Value = '''165.20'
$ws.Range("E36").Value = '  +1.78%  '
$ws.Range("D37").Value = '''19.92'
$ws.Range("E37").Value = '  +4.22%  '
$ws.Range("D38").Value = '''0.118'
$ws.Range("E38").Value = '  +0.69%  '
$ws.Range("D39").Value = '''19.38'
$ws.Range("E39").Value = '  +0.38%  '
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").Value = '''183.61'
$ws.Range("E40").Value = '  +11.91%  '
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").Value = '''0.343'
$ws.Range("E42").Value = '  +5.41%  '
$ws.Range("D43").Value = '''5.02'
$ws.Range("E43").Value = '  +2.75%  '
$ws.Range("D44").Value = '''1.67'
$ws.Range("E44").Value = '  +0.84%  '
$ws.Range("D45").Value = '''1.21'
$ws.Range("E45").Value = '  +3.25%  '
$ws.Range("D46").Value = '''39.99'
$ws.Range("E46").Value = '  +2.61%  '
$ws.Range("D47").Value = '''0.0871'
$ws.Range("E47").Value = '  +2.78%  '
$ws.Range("D48").Value = '''2.35'
$ws.Range("E48").Value = '  +0.04%  '
$ws.Range("D49").Value = '''0.573'
$ws.Range("E49").Value = '  +9.15%  '
$ws.Range("D50").Value = '''3.74'
$ws.Range("E50").Value = '  +4.01%  '
$ws.Range("D51").Value = '''0.644'
$ws.Range("E51").Value = '  +9.31%  '
